# Fixed button hover bug for safari
#
# - "Resume hover bug (Safari)" (row 12) is now resolved: style flips from
#   "Neutral" (in-progress, orange) to "Good" (done, green) and the
#   Date Resolved cell (F12) is filled in.
# - "Grammer and spelling check" (row 19) moves from "not started" (no
#   fill) to "Neutral" (in-progress, orange) as the next item being worked.
# - The two stale =TODAY() lookups (F4/F10) are frozen to their last
#   calculated value instead of staying volatile.
# - Selection cursor left where the editor last clicked.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Freeze the two volatile TODAY() cached cells to plain static values.
$ws.Range("F4").Value = 44321
$ws.Range("F10").Value = 44321

# Row 12 "Resume hover bug (Safari)": mark complete -> "Good" style, and
# fill in the Date Resolved column.
$row12 = $ws.Range("B12:F12")
$row12.Font.Color = 24832
$row12.Interior.Color = 13561798

# F12 needs the same date-formatted variant of the style that E12 already
# carries, so clone E12's format onto it before writing the resolved date.
$ws.Range("E12").Copy()
$ws.Range("F12").PasteSpecial(-4122)
$ws.Range("F12").Value = 44322

# Row 19 "Grammer and spelling check": move to "Neutral" (in progress)
# style.
$row19 = $ws.Range("B19:F19")
$row19.Font.Color = 22428
$row19.Interior.Color = 10284031

# Leave the selection where the author's last click landed.
$null = $ws.Range("I24").Select()
